$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing values (row 252-255, column C)
$ws.Range("C252").Value = 332052.35
$ws.Range("C253").Value = 91747.38
$ws.Range("C254").Value = 216021.3
$ws.Range("C255").Value = 102229.33

# Add new rows 256 and 257
$ws.Range("A256").Value = 6
$ws.Range("B256").Value = 2
$ws.Range("C256").Value = 3340
$ws.Range("D256").Value = 2025
$ws.Range("E256").Value = "Bibi Cell Manauara"

$ws.Range("A257").Value = 6
$ws.Range("B257").Value = 4
$ws.Range("C257").Value = 1800.01
$ws.Range("D257").Value = 2025
$ws.Range("E257").Value = "Bibi Cell Ponta Negra"
